$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.676.94'
$ws.Range("E2").Value = '  +0.66%  '

$ws.Range("D3").Value = '1.850.67'
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.60'
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4258'
$ws.Range("E7").Value = '  +0.56%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3632'
$ws.Range("E8").Value = '  +0.47%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '44.76'
$ws.Range("E9").Value = '  +2.72%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07294'
$ws.Range("E10").Value = '  +1.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8734'
$ws.Range("E11").Value = '  -2.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.60'
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").Value = '1.903.62'
$ws.Range("E13").Value = '  +5.03%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.311'
$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.509'
$ws.Range("E15").Value = '  -0.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06895'
$ws.Range("E16").Value = '  +1.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '79.74'
$ws.Range("E18").Value = '  +3.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009019'
$ws.Range("E19").Value = '  +1.62%  '

$ws.Range("E20").Value = '  +0.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.35'
$ws.Range("E21").Value = '  +0.37%  '

$ws.Range("D22").Value = '27.691.83'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.969'
$ws.Range("E23").Value = '  +1.07%  '

$ws.Range("E24").Value = '  -3.75%  '

$ws.Range("D25").Value = '2.071.24'
$ws.Range("E25").Value = '  +1.88%  '

$ws.Range("E26").Value = '  -3.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '153.84'
$ws.Range("E27").Value = '  +2.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.86'
$ws.Range("E28").Value = '  +3.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '121.72'
$ws.Range("E29").Value = '  +10.02%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.264'
$ws.Range("E30").Value = '  -0.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.869'
$ws.Range("E31").Value = '  +10.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08909'
$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7575'
$ws.Range("E33").Value = '  -1.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.967'
$ws.Range("E34").Value = '  +3.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.517'
$ws.Range("E35").Value = '  +1.03%  '

$ws.Range("E36").Value = '  +2.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05377'
$ws.Range("E37").Value = '  +0.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.089'
$ws.Range("E38").Value = '  -0.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01929'
$ws.Range("E39").Value = '  +0.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.818'
$ws.Range("E40").Value = '  -4.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5060'
$ws.Range("E41").Value = '  +0.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1648'
$ws.Range("E42").Value = '  +1.06%  '

$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.313'
$ws.Range("E44").Value = '  +1.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.06543'
$ws.Range("E45").Value = '  -1.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.35'
$ws.Range("E46").Value = '  +2.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '104.97'
$ws.Range("E47").Value = '  -0.56%  '

$ws.Range("B48").Value = 'Decentraland'
$ws.Range("C48").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4645'
$ws.Range("E48").Value = '  -1.25%  '

$ws.Range("B49").Value = 'PaxDollar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  +0.07%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.617'
$ws.Range("E50").Value = '  -1.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.26'
$ws.Range("E51").Value = '  +0.11%  '
